# ingredients.xlsx update: add new ingredients, rename a couple of
# existing ones, tweak column width and a few cell styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename ingredients -------------------------------------------------

# Chocolate: clarify the cocoa percentage on the sugar-free bar.
$ws.Range("B8").Value = "Amul Dark Chocolate 55% Sugar Free"

# Row 10 was mislabeled "Stabilizer" -- Weikfield Cornstarch is really a
# thickener, so relabel the ingredient category.
$ws.Range("A10").Value = "Thickener"

# --- 2. Column B is wider now (longer ingredient names) --------------------

$ws.Columns.Item(2).ColumnWidth = 32.25

# --- 3. Match formatting of row 8 / row 14 to the rest of the new entries --

$ws.Range("B8").HorizontalAlignment = -4108
$ws.Range("B8").VerticalAlignment = -4160

$ws.Range("A14:B14").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)

# --- 4. Append three new ingredient rows ------------------------------------

# Use row 13's formatting (borders/fonts/number formats) as the template for
# the new data rows, then overwrite values/formulas.
$ws.Range("A13:O13").Copy()
$ws.Range("A15:O15").PasteSpecial(-4122)
$ws.Range("A16:O16").PasteSpecial(-4122)
$ws.Range("A17:O17").PasteSpecial(-4122)

# Row 15: Amul Dark Chocolate 75%
$ws.Range("A15").Value = "Chocolate"
$ws.Range("B15").Value = "Amul Dark Chocolate 75% "
$ws.Range("C15").Value = 538.0
$ws.Range("D15").Value = 9.6
$ws.Range("E15").Value = 37.4
$ws.Range("F15").Value = 25.3
$ws.Range("G15").Value = 39.8
$ws.Range("H15").Value = 0.0
$ws.Range("I15").Value = 57.0
$ws.Range("J15").Value = 0.0
$ws.Range("K15").Value = 0.0
$ws.Range("L15").Value = 0.0
$ws.Range("M15").Formula = '=0.9*((24.5*100)+(0*100)+(0.8*190))/25.3'
$ws.Range("N15").Formula = '=0.95*((24.5*100)+(0*16)+(0.8*130))/25.3'
$ws.Range("O15").Value = 120.0
$ws.Range("M15:N15").NumberFormat = "0.00"

# Row 16: Amul Belgian Chocolate
$ws.Range("A16").Value = "Chocolate"
$ws.Range("B16").Value = "Amul Belgian Chocolate"
$ws.Range("C16").Value = 551.0
$ws.Range("D16").Value = 8.7
$ws.Range("E16").Value = 51.2
$ws.Range("F16").Value = 49.6
$ws.Range("G16").Value = 35.3
$ws.Range("H16").Value = 0.0
$ws.Range("I16").Value = 155.0
$ws.Range("J16").Value = 0.0
$ws.Range("K16").Formula = '=1.07*(0.88*D16+8)'
$ws.Range("L16").Value = 0.0
$ws.Range("M16").Formula = '=0.9*((40*100)+(8*100) + (1.6*190))/49.6'
$ws.Range("N16").Formula = '=0.95*((40*100)+(8*16)+(1.6*130))/49.6'
$ws.Range("O16").Value = 128.0
$ws.Range("K16").NumberFormat = "0.00"
$ws.Range("M16").NumberFormat = "0.00"
$ws.Range("N16").NumberFormat = "0"

# Row 17: Hershey's Cocoa Powder
$ws.Range("A17").Value = "Cocoa Powder"
$ws.Range("B17").Value = "Hershey's Cocoa Powder"
$ws.Range("C17").Value = 409.0
$ws.Range("D17").Value = 19.6
$ws.Range("E17").Value = 58.0
$ws.Range("F17").Value = 1.8
$ws.Range("G17").Value = 11.0
$ws.Range("H17").Value = 0.0
$ws.Range("I17").Value = 4.0
$ws.Range("J17").Value = 0.0
$ws.Range("K17").Value = 0.0
$ws.Range("L17").Value = 0.0
$ws.Range("M17").Value = 190.0
$ws.Range("N17").Value = 125.0
$ws.Range("O17").Value = 200.0
$ws.Range("K17").NumberFormat = "0.00"
$ws.Range("M17").NumberFormat = "0.00"
$ws.Range("N17").NumberFormat = "0"
